$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Bmp6"
$ws.Cells.Item(2,3).Value = "Bmpr1a"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 21.443826
$ws.Cells.Item(2,8).Value = 64.331478
$ws.Cells.Item(2,9).Value = 0.6062978927103765
$ws.Cells.Item(2,10).Value = 0.6062978927103765
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.330485333333333
$ws.Cells.Item(2,14).Value = 6.991456
$ws.Cells.Item(2,15).Value = 0.03655621216383393
$ws.Cells.Item(2,16).Value = 0.03655621216383393
$ws.Cells.Item(2,17).Value = 49.974521983552
$ws.Cells.Item(2,18).Value = 449.7706978519681
$ws.Cells.Item(2,19).Value = 0.02216395440040594
$ws.Cells.Item(2,20).Value = 0.02216395440040594

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Bmp6"
$ws.Cells.Item(3,3).Value = "Bmpr1a"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 21.443826
$ws.Cells.Item(3,8).Value = 64.331478
$ws.Cells.Item(3,9).Value = 0.6062978927103765
$ws.Cells.Item(3,10).Value = 0.6062978927103765
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 29.178763
$ws.Cells.Item(3,14).Value = 87.53628900000001
$ws.Cells.Item(3,15).Value = 0.4577008212193115
$ws.Cells.Item(3,16).Value = 0.4577008212193115
$ws.Cells.Item(3,17).Value = 625.7043166672381
$ws.Cells.Item(3,18).Value = 5631.338850005143
$ws.Cells.Item(3,19).Value = 0.2775030433970773
$ws.Cells.Item(3,20).Value = 0.2775030433970773

# Row 4: ECs -> M1
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Bmp6"
$ws.Cells.Item(4,3).Value = "Bmpr1a"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 21.443826
$ws.Cells.Item(4,8).Value = 64.331478
$ws.Cells.Item(4,9).Value = 0.6062978927103765
$ws.Cells.Item(4,10).Value = 0.6062978927103765
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.3937893333333333
$ws.Cells.Item(4,14).Value = 1.181368
$ws.Cells.Item(4,15).Value = 0.006177016525822971
$ws.Cells.Item(4,16).Value = 0.006177016525822971
$ws.Cells.Item(4,17).Value = 8.444349944656
$ws.Cells.Item(4,18).Value = 75.999149501904
$ws.Cells.Item(4,19).Value = 0.003745112102843638
$ws.Cells.Item(4,20).Value = 0.003745112102843638

# Row 5: ECs -> M2
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Bmp6"
$ws.Cells.Item(5,3).Value = "Bmpr1a"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 21.443826
$ws.Cells.Item(5,8).Value = 64.331478
$ws.Cells.Item(5,9).Value = 0.6062978927103765
$ws.Cells.Item(5,10).Value = 0.6062978927103765
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.4229803333333333
$ws.Cells.Item(5,14).Value = 1.268941
$ws.Cells.Item(5,15).Value = 0.006634909297775398
$ws.Cells.Item(5,16).Value = 0.006634909297775398
$ws.Cells.Item(5,17).Value = 9.070316669422001
$ws.Cells.Item(5,18).Value = 81.63285002479802
$ws.Cells.Item(5,19).Value = 0.004022731525565707
$ws.Cells.Item(5,20).Value = 0.004022731525565707

# Row 6: ECs -> Neutro
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Bmp6"
$ws.Cells.Item(6,3).Value = "Bmpr1a"
$ws.Cells.Item(6,4).Value = "Neutro"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 21.443826
$ws.Cells.Item(6,8).Value = 64.331478
$ws.Cells.Item(6,9).Value = 0.6062978927103765
$ws.Cells.Item(6,10).Value = 0.6062978927103765
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 9.749056666666668
$ws.Cells.Item(6,14).Value = 29.24717
$ws.Cells.Item(6,15).Value = 0.1529246199520842
$ws.Cells.Item(6,16).Value = 0.1529246199520842
$ws.Cells.Item(6,17).Value = 209.05707482414
$ws.Cells.Item(6,18).Value = 1881.51367341726
$ws.Cells.Item(6,19).Value = 0.09271787482048387
$ws.Cells.Item(6,20).Value = 0.09271787482048387

# Row 7: ECs -> sCs
$ws.Cells.Item(7,1).Value = "ECs"
$ws.Cells.Item(7,2).Value = "Bmp6"
$ws.Cells.Item(7,3).Value = "Bmpr1a"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 21.443826
$ws.Cells.Item(7,8).Value = 64.331478
$ws.Cells.Item(7,9).Value = 0.6062978927103765
$ws.Cells.Item(7,10).Value = 0.6062978927103765
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 21.67565866666666
$ws.Cells.Item(7,14).Value = 65.02697599999999
$ws.Cells.Item(7,15).Value = 0.3400064208411719
$ws.Cells.Item(7,16).Value = 0.3400064208411719
$ws.Cells.Item(7,17).Value = 464.809052883392
$ws.Cells.Item(7,18).Value = 4183.281475950527
$ws.Cells.Item(7,19).Value = 0.206145176464
$ws.Cells.Item(7,20).Value = 0.206145176464

# Row 8: FAPs -> ECs
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Bmp6"
$ws.Cells.Item(8,3).Value = "Bmpr1a"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 3.641794
$ws.Cells.Item(8,8).Value = 10.925382
$ws.Cells.Item(8,9).Value = 0.1029672609675761
$ws.Cells.Item(8,10).Value = 0.1029672609675761
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.330485333333333
$ws.Cells.Item(8,14).Value = 6.991456
$ws.Cells.Item(8,15).Value = 0.03655621216383393
$ws.Cells.Item(8,16).Value = 0.03655621216383393
$ws.Cells.Item(8,17).Value = 8.487147504021333
$ws.Cells.Item(8,18).Value = 76.384327536192
$ws.Cells.Item(8,19).Value = 0.003764093037859567
$ws.Cells.Item(8,20).Value = 0.003764093037859567

# Row 9: FAPs -> FAPs
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Bmp6"
$ws.Cells.Item(9,3).Value = "Bmpr1a"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 3.641794
$ws.Cells.Item(9,8).Value = 10.925382
$ws.Cells.Item(9,9).Value = 0.1029672609675761
$ws.Cells.Item(9,10).Value = 0.1029672609675761
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 29.178763
$ws.Cells.Item(9,14).Value = 87.53628900000001
$ws.Cells.Item(9,15).Value = 0.4577008212193115
$ws.Cells.Item(9,16).Value = 0.4577008212193115
$ws.Cells.Item(9,17).Value = 106.263044020822
$ws.Cells.Item(9,18).Value = 956.3673961873982
$ws.Cells.Item(9,19).Value = 0.04712819990356272
$ws.Cells.Item(9,20).Value = 0.04712819990356272

# Row 10: FAPs -> M1
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Bmp6"
$ws.Cells.Item(10,3).Value = "Bmpr1a"
$ws.Cells.Item(10,4).Value = "M1"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 3.641794
$ws.Cells.Item(10,8).Value = 10.925382
$ws.Cells.Item(10,9).Value = 0.1029672609675761
$ws.Cells.Item(10,10).Value = 0.1029672609675761
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.3937893333333333
$ws.Cells.Item(10,14).Value = 1.181368
$ws.Cells.Item(10,15).Value = 0.006177016525822971
$ws.Cells.Item(10,16).Value = 0.006177016525822971
$ws.Cells.Item(10,17).Value = 1.434099631397334
$ws.Cells.Item(10,18).Value = 12.906896682576
$ws.Cells.Item(10,19).Value = 0.0006360304726154438
$ws.Cells.Item(10,20).Value = 0.0006360304726154438

# Row 11: FAPs -> M2
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Bmp6"
$ws.Cells.Item(11,3).Value = "Bmpr1a"
$ws.Cells.Item(11,4).Value = "M2"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 3.641794
$ws.Cells.Item(11,8).Value = 10.925382
$ws.Cells.Item(11,9).Value = 0.1029672609675761
$ws.Cells.Item(11,10).Value = 0.1029672609675761
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.4229803333333333
$ws.Cells.Item(11,14).Value = 1.268941
$ws.Cells.Item(11,15).Value = 0.006634909297775398
$ws.Cells.Item(11,16).Value = 0.006634909297775398
$ws.Cells.Item(11,17).Value = 1.540407240051334
$ws.Cells.Item(11,18).Value = 13.863665160462
$ws.Cells.Item(11,19).Value = 0.0006831784371602363
$ws.Cells.Item(11,20).Value = 0.0006831784371602363

# Row 12: FAPs -> Neutro
$ws.Cells.Item(12,1).Value = "FAPs"
$ws.Cells.Item(12,2).Value = "Bmp6"
$ws.Cells.Item(12,3).Value = "Bmpr1a"
$ws.Cells.Item(12,4).Value = "Neutro"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 3.641794
$ws.Cells.Item(12,8).Value = 10.925382
$ws.Cells.Item(12,9).Value = 0.1029672609675761
$ws.Cells.Item(12,10).Value = 0.1029672609675761
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 9.749056666666668
$ws.Cells.Item(12,14).Value = 29.24717
$ws.Cells.Item(12,15).Value = 0.1529246199520842
$ws.Cells.Item(12,16).Value = 0.1529246199520842
$ws.Cells.Item(12,17).Value = 35.50405607432668
$ws.Cells.Item(12,18).Value = 319.5365046689401
$ws.Cells.Item(12,19).Value = 0.01574622925097365
$ws.Cells.Item(12,20).Value = 0.01574622925097365

# Row 13: FAPs -> sCs
$ws.Cells.Item(13,1).Value = "FAPs"
$ws.Cells.Item(13,2).Value = "Bmp6"
$ws.Cells.Item(13,3).Value = "Bmpr1a"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 3.641794
$ws.Cells.Item(13,8).Value = 10.925382
$ws.Cells.Item(13,9).Value = 0.1029672609675761
$ws.Cells.Item(13,10).Value = 0.1029672609675761
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 21.67565866666666
$ws.Cells.Item(13,14).Value = 65.02697599999999
$ws.Cells.Item(13,15).Value = 0.3400064208411719
$ws.Cells.Item(13,16).Value = 0.3400064208411719
$ws.Cells.Item(13,17).Value = 78.93828367831466
$ws.Cells.Item(13,18).Value = 710.4445531048319
$ws.Cells.Item(13,19).Value = 0.03500952986540445
$ws.Cells.Item(13,20).Value = 0.03500952986540445

# Row 14: sCs -> ECs
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Bmp6"
$ws.Cells.Item(14,3).Value = "Bmpr1a"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 10.28284533333333
$ws.Cells.Item(14,8).Value = 30.848536
$ws.Cells.Item(14,9).Value = 0.2907348463220475
$ws.Cells.Item(14,10).Value = 0.2907348463220475
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 2.330485333333333
$ws.Cells.Item(14,14).Value = 6.991456
$ws.Cells.Item(14,15).Value = 0.03655621216383393
$ws.Cells.Item(14,16).Value = 0.03655621216383393
$ws.Cells.Item(14,17).Value = 23.96402023426844
$ws.Cells.Item(14,18).Value = 215.676182108416
$ws.Cells.Item(14,19).Value = 0.01062816472556842
$ws.Cells.Item(14,20).Value = 0.01062816472556842

# Row 15: sCs -> FAPs
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Bmp6"
$ws.Cells.Item(15,3).Value = "Bmpr1a"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 10.28284533333333
$ws.Cells.Item(15,8).Value = 30.848536
$ws.Cells.Item(15,9).Value = 0.2907348463220475
$ws.Cells.Item(15,10).Value = 0.2907348463220475
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 29.178763
$ws.Cells.Item(15,14).Value = 87.53628900000001
$ws.Cells.Item(15,15).Value = 0.4577008212193115
$ws.Cells.Item(15,16).Value = 0.4577008212193115
$ws.Cells.Item(15,17).Value = 300.0407069469894
$ws.Cells.Item(15,18).Value = 2700.366362522904
$ws.Cells.Item(15,19).Value = 0.1330695779186715
$ws.Cells.Item(15,20).Value = 0.1330695779186715

# Row 16: sCs -> M1
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Bmp6"
$ws.Cells.Item(16,3).Value = "Bmpr1a"
$ws.Cells.Item(16,4).Value = "M1"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 10.28284533333333
$ws.Cells.Item(16,8).Value = 30.848536
$ws.Cells.Item(16,9).Value = 0.2907348463220475
$ws.Cells.Item(16,10).Value = 0.2907348463220475
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.3937893333333333
$ws.Cells.Item(16,14).Value = 1.181368
$ws.Cells.Item(16,15).Value = 0.006177016525822971
$ws.Cells.Item(16,16).Value = 0.006177016525822971
$ws.Cells.Item(16,17).Value = 4.04927480858311
$ws.Cells.Item(16,18).Value = 36.443473277248
$ws.Cells.Item(16,19).Value = 0.001795873950363889
$ws.Cells.Item(16,20).Value = 0.001795873950363889

# Row 17: sCs -> M2
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Bmp6"
$ws.Cells.Item(17,3).Value = "Bmpr1a"
$ws.Cells.Item(17,4).Value = "M2"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 10.28284533333333
$ws.Cells.Item(17,8).Value = 30.848536
$ws.Cells.Item(17,9).Value = 0.2907348463220475
$ws.Cells.Item(17,10).Value = 0.2907348463220475
$ws.Cells.Item(17,11).Value = 1
$ws.Cells.Item(17,12).Value = 0.3333333333333333
$ws.Cells.Item(17,13).Value = 0.4229803333333333
$ws.Cells.Item(17,14).Value = 1.268941
$ws.Cells.Item(17,15).Value = 0.006634909297775398
$ws.Cells.Item(17,16).Value = 0.006634909297775398
$ws.Cells.Item(17,17).Value = 4.349441346708444
$ws.Cells.Item(17,18).Value = 39.144972120376
$ws.Cells.Item(17,19).Value = 0.001928999335049455
$ws.Cells.Item(17,20).Value = 0.001928999335049455

# Row 18: sCs -> Neutro
$ws.Cells.Item(18,1).Value = "sCs"
$ws.Cells.Item(18,2).Value = "Bmp6"
$ws.Cells.Item(18,3).Value = "Bmpr1a"
$ws.Cells.Item(18,4).Value = "Neutro"
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 10.28284533333333
$ws.Cells.Item(18,8).Value = 30.848536
$ws.Cells.Item(18,9).Value = 0.2907348463220475
$ws.Cells.Item(18,10).Value = 0.2907348463220475
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 9.749056666666668
$ws.Cells.Item(18,14).Value = 29.24717
$ws.Cells.Item(18,15).Value = 0.1529246199520842
$ws.Cells.Item(18,16).Value = 0.1529246199520842
$ws.Cells.Item(18,17).Value = 100.2480418492356
$ws.Cells.Item(18,18).Value = 902.2323766431201
$ws.Cells.Item(18,19).Value = 0.04446051588062674
$ws.Cells.Item(18,20).Value = 0.04446051588062674

# Row 19: sCs -> sCs
$ws.Cells.Item(19,1).Value = "sCs"
$ws.Cells.Item(19,2).Value = "Bmp6"
$ws.Cells.Item(19,3).Value = "Bmpr1a"
$ws.Cells.Item(19,4).Value = "sCs"
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 10.28284533333333
$ws.Cells.Item(19,8).Value = 30.848536
$ws.Cells.Item(19,9).Value = 0.2907348463220475
$ws.Cells.Item(19,10).Value = 0.2907348463220475
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 21.67565866666666
$ws.Cells.Item(19,14).Value = 65.02697599999999
$ws.Cells.Item(19,15).Value = 0.3400064208411719
$ws.Cells.Item(19,16).Value = 0.3400064208411719
$ws.Cells.Item(19,17).Value = 222.8874455674595
$ws.Cells.Item(19,18).Value = 2005.987010107136
$ws.Cells.Item(19,19).Value = 0.09885171451176755
$ws.Cells.Item(19,20).Value = 0.09885171451176755
